# "updated confusion matrices and accuracies"
# Rows 4-7 (D:K) hold the raw accuracy fractions used throughout the sheet;
# rows 16-19 (the "Percentages" block) already recompute from these via
# shared formulas (e.g. D16 = D4*100), so only the raw values below need
# to be written and the percentage rows will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    4 = @{ D = 0.977621483375959;  E = 0.85741687979539605; F = 0.906649616368286;
           G = 0.92838874680306904; H = 0.96419437340153402; I = 0.89897698209718602;
           J = 0.56329923273657201; K = 0.64002557544756999 }
    5 = @{ D = 0.90025575447570305; E = 0.50447570332480796; F = 0.74360613810741605;
           G = 0.81265984654731405; H = 0.83184143222506302; I = 0.83312020460358005;
           J = 0.37020460358056201 }
    6 = @{ D = 0.92199488491048598; E = 0.68734015345268495; F = 0.83759590792838801;
           G = 0.87851662404092001; H = 0.91432225063938599; I = 0.90345268542199397;
           J = 0.62020460358056195 }
    7 = @{ D = 0.98399999999999999; E = 0.873142857142857;   F = 0.91314285714285703;
           G = 0.94399999999999995; H = 0.94628571428571395; I = 0.91542857142857104;
           J = 0.54285714285714204 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value2 = $newValues[$row][$col]
    }
}

# The author's active selection moved from D21 to H10 before saving.
$ws.Range("H10").Select()
